$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.202777504920959
$ws.Range("B1").Value = 1.887801051139832
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.857860088348389
$ws.Range("E1").Value = 1.20584511756897
